$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, C, D (rows 2-25)
$arrBD = New-Object 'object[,]' 24,3
$arrBD[0,0] = 1.213056729844538; $arrBD[0,1] = 0.3228198103155648; $arrBD[0,2] = 0.01941460295053332
$arrBD[1,0] = 1.065510157487665; $arrBD[1,1] = 0.2833549914221578; $arrBD[1,2] = 0.01736633054169801
$arrBD[2,0] = 0.9747452631027613; $arrBD[2,1] = 0.259017211634017; $arrBD[2,2] = 0.01610683451141171
$arrBD[3,0] = 0.9377162116397812; $arrBD[3,1] = 0.2490731027117192; $arrBD[3,2] = 0.01559314843509441
$arrBD[4,0] = 0.9315650912813567; $arrBD[4,1] = 0.2474203178887819; $arrBD[4,2] = 0.01550782621929869
$arrBD[5,0] = 0.9742460431247082; $arrBD[5,1] = 0.2588832076546908; $arrBD[5,2] = 0.01609990846147724
$arrBD[6,0] = 1.162218828581103; $arrBD[6,1] = 0.3092346387274176; $arrBD[6,2] = 0.0187087614711885
$arrBD[7,0] = 1.529440664866854; $arrBD[7,1] = 0.4071164138300674; $arrBD[7,2] = 0.02380887256085629
$arrBD[8,0] = 1.79836787514472; $arrBD[8,1] = 0.4784957298887775; $arrBD[8,2] = 0.02754513163202432
$arrBD[9,0] = 1.920518466643898; $arrBD[9,1] = 0.5108501053937289; $arrBD[9,2] = 0.02924231419535062
$arrBD[10,0] = 1.96674624670851; $arrBD[10,1] = 0.5230848207709187; $arrBD[10,2] = 0.02988461615090898
$arrBD[11,0] = 1.956791532486648; $arrBD[11,1] = 0.5204506271037985; $arrBD[11,2] = 0.02974630249445198
$arrBD[12,0] = 1.924322222300134; $arrBD[12,1] = 0.5118570096828421; $arrBD[12,2] = 0.02929516463519377
$arrBD[13,0] = 1.904430148518429; $arrBD[13,1] = 0.5065909200399119; $arrBD[13,2] = 0.02901877874589331
$arrBD[14,0] = 1.790381198118268; $arrBD[14,1] = 0.4763789144324164; $arrBD[14,2] = 0.02743416462142534
$arrBD[15,0] = 1.720367359308625; $arrBD[15,1] = 0.4578146885784804; $arrBD[15,2] = 0.02646140297542132
$arrBD[16,0] = 1.680079842492546; $arrBD[16,1] = 0.4471261033194196; $arrBD[16,2] = 0.02590166636442603
$arrBD[17,0] = 1.666436231437274; $arrBD[17,1] = 0.4435052677839053; $arrBD[17,2] = 0.02571211063713008
$arrBD[18,0] = 1.727822263141434; $arrBD[18,1] = 0.4597920169316012; $arrBD[18,2] = 0.026564979126114
$arrBD[19,0] = 1.933860015709854; $arrBD[19,1] = 0.5143816332765709; $arrBD[19,2] = 0.02942768542513363
$arrBD[20,0] = 2.068353832317882; $arrBD[20,1] = 0.5499587420608805; $arrBD[20,2] = 0.03129637649145423
$arrBD[21,0] = 1.996587351909056; $arrBD[21,1] = 0.5309798956201348; $arrBD[21,2] = 0.03029923680481517
$arrBD[22,0] = 1.724452012012534; $arrBD[22,1] = 0.4588981158648267; $arrBD[22,2] = 0.02651815385038958
$arrBD[23,0] = 1.43024888320673; $arrBD[23,1] = 0.3807298380316411; $arrBD[23,2] = 0.02243097516476666
$ws.Range("B2:D25").Value = $arrBD

# Columns F, G, H (rows 2-25)
$arrFH = New-Object 'object[,]' 24,3
$arrFH[0,0] = 0.430388737815619; $arrFH[0,1] = 0.2777800091919644; $arrFH[0,2] = 0.4417816666615906
$arrFH[1,0] = 0.4328886863516956; $arrFH[1,1] = 0.2810209571142295; $arrFH[1,2] = 0.4500004064366578
$arrFH[2,0] = 0.4350970362268711; $arrFH[2,1] = 0.2836156918042221; $arrFH[2,2] = 0.4555491057942831
$arrFH[3,0] = 0.4361652221734147; $arrFH[3,1] = 0.2848238556574429; $arrFH[3,2] = 0.4579360904076637
$arrFH[4,0] = 0.4363527253908401; $arrFH[4,1] = 0.2850335398112023; $arrFH[4,2] = 0.4583400340787236
$arrFH[5,0] = 0.4351107622628376; $arrFH[5,1] = 0.2836313766978051; $arrFH[5,2] = 0.4555807886008481
$arrFH[6,0] = 0.4311103241399579; $arrFH[6,1] = 0.2787711914752933; $arrFH[6,2] = 0.4445109365135664
$arrFH[7,0] = 0.4286565667544195; $arrFH[7,1] = 0.2740971149103473; $arrFH[7,2] = 0.4268103748463687
$arrFH[8,0] = 0.4302076403426653; $arrFH[8,1] = 0.2737042005166188; $arrFH[8,2] = 0.4162781369516324
$arrFH[9,0] = 0.431656005885003; $arrFH[9,1] = 0.2742030514762632; $arrFH[9,2] = 0.4120299100306255
$arrFH[10,0] = 0.4323124325517824; $arrFH[10,1] = 0.2744907836870425; $arrFH[10,2] = 0.4104998135478723
$arrFH[11,0] = 0.4321662393469836; $arrFH[11,1] = 0.2744243978458485; $arrFH[11,2] = 0.4108258421714481
$arrFH[12,0] = 0.4317078401858296; $arrFH[12,1] = 0.2742247347462268; $arrFH[12,2] = 0.4119024494152086
$arrFH[13,0] = 0.4314411514865171; $arrFH[13,1] = 0.2741153469670081; $arrFH[13,2] = 0.4125721564648046
$arrFH[14,0] = 0.4301280312677136; $arrFH[14,1] = 0.2736853532468899; $arrFH[14,2] = 0.4165667401346838
$arrFH[15,0] = 0.4295135368886136; $arrFH[15,1] = 0.2735961034271526; $arrFH[15,2] = 0.4191567325595287
$arrFH[16,0] = 0.429229940611819; $arrFH[16,1] = 0.273608464108321; $arrFH[16,2] = 0.4206974995174591
$arrFH[17,0] = 0.4291458805542732; $arrFH[17,1] = 0.2736235464545302; $arrFH[17,2] = 0.4212279319078078
$arrFH[18,0] = 0.4295717138565607; $arrFH[18,1] = 0.2735990019647545; $arrFH[18,2] = 0.4188757337847733
$arrFH[19,0] = 0.4318395438949452; $arrFH[19,1] = 0.274280687383154; $arrFH[19,2] = 0.4115840858751483
$arrFH[20,0] = 0.4339515215159366; $arrFH[20,1] = 0.2753028998018721; $arrFH[20,2] = 0.4072770512113806
$arrFH[21,0] = 0.4327663103649755; $arrFH[21,1] = 0.2747040988571285; $arrFH[21,2] = 0.4095336670924468
$arrFH[22,0] = 0.4295451950711424; $arrFH[22,1] = 0.2735974933082588; $arrFH[22,2] = 0.4190026121808046
$arrFH[23,0] = 0.4287362588925561; $arrFH[23,1] = 0.2748333215114087; $arrFH[23,2] = 0.4311667942157271
$ws.Range("F2:H25").Value = $arrFH

Write-Host "Updated pl_mw line-loading values for the 380 kV case"
